# Re-sort the results table ("Table1") by the "Факултетен номер" column
# (column B) ascending, instead of the original sort by "Име" (column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects("Table1")

$table.Sort.SortFields.Clear()
$table.Sort.SortFields.Add($ws.Range("B1:B20"), 0, 1) | Out-Null
$table.Sort.Header = 1
$table.Sort.Apply()

$ws.Range("A1").Select() | Out-Null
